# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.129.62"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.547.04"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.83"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.57"
$ws.Range("E6").Value = "  -4.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.546.09"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("E10").Value = "  -2.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.09"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.150.26"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000183"
$ws.Range("E14").Value = "  -3.13%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.81"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.549.15"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.250.48"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.95"
$ws.Range("E19").Value = "  -3.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.35"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.81"
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.07"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.577"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.690.79"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.98"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.73"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.54"
$ws.Range("E29").Value = "  +23.74%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.50"
$ws.Range("E31").Value = "  +2.22%  "
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.550.83"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.00"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "169.91"
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.90"
$ws.Range("E38").Value = "  -2.29%  "
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.00"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.49"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.24"
$ws.Range("E45").Value = "  +3.14%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.64"
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.457.65"
$ws.Range("E49").Value = "  +3.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.90"
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("E51").Value = "  +0.56%  "
